$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "177÷6=" "676÷7="
Replace-Text "753÷9=" "863÷4="
Replace-Text "239÷8=" "563÷9="
Replace-Text "616÷7=" "146÷4="
Replace-Text "806÷9=" "756÷8="
Replace-Text "925÷6=" "520÷4="
Replace-Text "810÷7=" "868÷3="
Replace-Text "168÷2=" "551÷5="
Replace-Text "185÷8=" "531÷5="
Replace-Text "448÷7=" "839÷8="
Replace-Text "433÷4=" "487÷8="
Replace-Text "748÷9=" "587÷6="
Replace-Text "313÷4=" "584÷8="
Replace-Text "815÷3=" "943÷5="
Replace-Text "801÷9=" "980÷8="
Replace-Text "194÷6=" "664÷4="
Replace-Text "918÷5=" "879÷8="
Replace-Text "759÷6=" "718÷3="
Replace-Text "483÷9=" "773÷5="
Replace-Text "808÷8=" "669÷7="
Replace-Text "691÷3=" "272÷7="
Replace-Text "946÷3=" "684÷4="
Replace-Text "140÷8=" "312÷5="
Replace-Text "194÷4=" "168÷4="
Replace-Text "726÷6=" "531÷5="
